$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 92, shifting existing rows 92-164 down to 93-165.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new weekly record.
$ws.Cells.Item(92, 1).Value = 5
$ws.Cells.Item(92, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(92, 3).Value = "Maule"
$ws.Cells.Item(92, 4).Value = 44447
$ws.Cells.Item(92, 5).Value = 7
$ws.Cells.Item(92, 6).Value = 100114014
$ws.Cells.Item(92, 7).Value = "Betarraga"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 4000
$ws.Cells.Item(92, 11).Value = 650
$ws.Cells.Item(92, 12).Value = 650
$ws.Cells.Item(92, 13).Value = 650
$ws.Cells.Item(92, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(92, 15).Value = "Región del Maule"
$ws.Cells.Item(92, 16).Value = 130
$ws.Cells.Item(92, 17).Value = 5
$ws.Cells.Item(92, 18).Value = "Hortaliza"

# Make sure the new date cell keeps the same date display format as the rest of column D.
$ws.Cells.Item(92, 4).NumberFormat = $ws.Cells.Item(93, 4).NumberFormat
